# Update "carjacking-by-neighborhood-by-month" workbook:
#  - refresh the "through <date>" labels (sheet name + header cell) from Sep 15 -> Sep 17
#  - update/add monthly carjacking counts for several neighborhoods (new data for 2022-09-25 pull)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / label updates -------------------------------------------------
$ws.Name = "Through 2022-09-17"
$ws.Range("B1").Value = "September 2022 (through September 17)"

# --- Updated counts for existing cells -------------------------------------
$updates = @{
    "B2"   = 4
    "K2"   = 7
    "T2"   = 5
    "K3"   = 10
    "AU3"  = 3
    "BD6"  = 2
    "B8"   = 6
    "K8"   = 4
    "AL8"  = 2
    "B10"  = 3
    "AC10" = 4
    "B15"  = 3
    "B17"  = 2
    "K18"  = 3
    "B22"  = 2
    "AC22" = 2
    "AC24" = 3
    "T29"  = 5
    "T33"  = 3
    "K37"  = 2
    "T50"  = 3
    "T57"  = 3
    "B89"  = 3
    "B96"  = 3
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Newly populated cells (previously blank) -------------------------------
$newValues = @{
    "BM4"  = 1
    "BM12" = 1
    "BD14" = 1
    "AU26" = 1
    "BD33" = 1
    "AL34" = 1
    "AL47" = 1
    "AC54" = 1
    "K55"  = 1
    "AU58" = 1
    "AC72" = 1
    "AC73" = 1
    "K77"  = 1
    "BD78" = 1
    "K83"  = 1
    "B85"  = 1
    "B90"  = 1
    "K17"  = 1
    "T97"  = 1
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
